$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string for the IN (India / B.1.617.2) variant column header
$ws.Range("G1").Value = "B.1.617.2"

# Row 2 - California: update existing counts and add new variant count
$ws.Range("B2").Value = 53
$ws.Range("D2").Value = 10
$ws.Range("G2").Value = 2

# Row 3 - Los Angeles: update existing counts and add new variant count
$ws.Range("B3").Value = 55
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 10
$ws.Range("G3").Value = 1

# Row 4 - San Diego: add new variant count
$ws.Range("G4").Value = 2

# Row 5 - San Francisco: add new variant count
$ws.Range("G5").Value = 1

# Move the active selection like the authored workbook shows
[void]$ws.Range("H4").Select()
